$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns before the existing "City" column (I),
# shifting City from I to L.
$ws.Columns("I:K").Insert()

# Header row
$ws.Range("I1").Value = "mood"
$ws.Range("J1").Value = "raise"
$ws.Range("K1").Value = "group"

# Data rows
$ws.Range("I2").Value = "happy"
$ws.Range("J2").Value = "no"
$ws.Range("K2").Value = "young_female"

$ws.Range("I3").Value = "happy"
$ws.Range("J3").Value = "yes"
$ws.Range("K3").Value = "mature_male"

$ws.Range("I4").Value = "happy"
$ws.Range("J4").Value = "yes"
$ws.Range("K4").Value = "young_male"

$ws.Range("I5").Value = "happy"
$ws.Range("J5").Value = "no"
$ws.Range("K5").Value = "mature_female"

$ws.Range("I6").Value = "happy"
$ws.Range("J6").Value = "yes"
$ws.Range("K6").Value = "young_female"

$ws.Range("I7").Value = "happy"
$ws.Range("J7").Value = "yes"

$ws.Range("I8").Value = "happy"
$ws.Range("J8").Value = "yes"
$ws.Range("K8").Value = "mature_female"

$ws.Range("I9").Value = "happy"
$ws.Range("K9").Value = "young_female"

$ws.Range("I10").Value = "happy"
$ws.Range("J10").Value = "yes"
$ws.Range("K10").Value = "young_male"

$ws.Range("I11").Value = "happy"
$ws.Range("J11").Value = "no"
$ws.Range("K11").Value = "mature_male"
